$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This report drops the tracking rows for the "2ef9fcc0-8841-4f9c-85f6-
# 6390a597d252" source file (it left the handback batch) on every sheet and
# refreshes the handoff/handback timestamps for the file that remains
# ("03f9a1cb-e8aa-4130-baf0-57c4d76a84e1").
# ---------------------------------------------------------------------------

# ---- Sheet 1: Overview ----------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Remove every hyperlink on the sheet (this engine's Hyperlinks.Delete()
# clears the whole collection no matter what range it is invoked from), we
# will re-create the ones that must survive once row 3 is gone.
$ws1.Hyperlinks.Delete()

# Drop row 3 (the 2ef9fcc0... entry); this naturally shrinks the dimension
# and reindexes/garbage-collects the shared strings table.
$ws1.Rows.Item(3).Delete()

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8daaac3f6d8c8c61bedb240fc113a5c5041c6a8b/e2e/03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.md", "", "", "03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.md") | Out-Null

# ---- Sheet 2: zh-cn ---------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()
$ws2.Rows.Item(3).Delete()

# Refresh the handoff/handback datetimes for the remaining row.
$ws2.Range("E2").Value = "2016-03-24 17:01:15"
$ws2.Range("H2").Value = "2016-03-24 17:02:01"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8daaac3f6d8c8c61bedb240fc113a5c5041c6a8b/e2e/03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.md", "", "", "03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/512b1536dba60c9c2ce1b99df8b50f6936b78c57/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.58559e9ba60de1b8adcfaef4d30383980f106c5c.zh-cn.xlf", "", "", "03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.58559e9ba60de1b8adcfaef4d30383980f106c5c.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/77b138b9e8a096290e23a320abe3311e6e001ba3/e2e/03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.md", "", "", "03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4b152161008677af0443f16b2ad89f0ef3496847/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.58559e9ba60de1b8adcfaef4d30383980f106c5c.zh-cn.xlf", "", "", "03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.58559e9ba60de1b8adcfaef4d30383980f106c5c.zh-cn.xlf") | Out-Null

# ---- Sheet 3: de-de ---------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()
$ws3.Rows.Item(3).Delete()

$ws3.Range("E2").Value = "2016-03-24 17:01:21"
$ws3.Range("H2").Value = "2016-03-24 17:02:13"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8daaac3f6d8c8c61bedb240fc113a5c5041c6a8b/e2e/03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.md", "", "", "03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/af971eabcc36ed1aa87cf202fb801c0437211839/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.58559e9ba60de1b8adcfaef4d30383980f106c5c.de-de.xlf", "", "", "03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.58559e9ba60de1b8adcfaef4d30383980f106c5c.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1e1671f8aa34588ac16ed2789ed10a18a3ef15c8/e2e/03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.md", "", "", "03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/907e138e0093d6217b16154460cbc8eba20fc592/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.58559e9ba60de1b8adcfaef4d30383980f106c5c.de-de.xlf", "", "", "03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.58559e9ba60de1b8adcfaef4d30383980f106c5c.de-de.xlf") | Out-Null

Write-Host "Done."
